# Append four new data rows (222-225) to Sheet1, following the same pattern
# as the existing rows: weekly Melón price records for
# "Terminal Hortofrutícola Agro Chillán".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ H = "Calameño"; I = "Primera"; J = 1000; K = 700; L = 800; M = 750; P = 750 },
    @{ H = "Calameño"; I = "Segunda"; J = 600;  K = 500; L = 600; M = 550; P = 550 },
    @{ H = "Tuna";     I = "Primera"; J = 1000; K = 700; L = 800; M = 750; P = 750 },
    @{ H = "Tuna";     I = "Segunda"; J = 600;  K = 500; L = 600; M = 550; P = 550 }
)

$startRow = 222

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]

    $ws.Cells.Item($r, 1).Value = 7
    $ws.Cells.Item($r, 2).Value = "Terminal Hortofrutícola Agro Chillán"
    $ws.Cells.Item($r, 3).Value = "Ñuble"

    $ws.Cells.Item($r, 4).Value = 44628
    $ws.Cells.Item($r, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"

    $ws.Cells.Item($r, 5).Value = 16
    $ws.Cells.Item($r, 6).Value = 100112027
    $ws.Cells.Item($r, 7).Value = "Melón"
    $ws.Cells.Item($r, 8).Value = $data.H
    $ws.Cells.Item($r, 9).Value = $data.I
    $ws.Cells.Item($r, 10).Value = $data.J
    $ws.Cells.Item($r, 11).Value = $data.K
    $ws.Cells.Item($r, 12).Value = $data.L
    $ws.Cells.Item($r, 13).Value = $data.M
    $ws.Cells.Item($r, 14).Value = "`$/unidad"
    $ws.Cells.Item($r, 15).Value = "Región del Maule"
    $ws.Cells.Item($r, 16).Value = $data.P
    $ws.Cells.Item($r, 17).Value = 1
    $ws.Cells.Item($r, 18).Value = "Hortaliza"
}
